# The sheet logs daily price observations for "Zapallo italiano" at the
# Feria Lagunitas de Puerto Montt. A new weekly observation (dated 45142,
# i.e. 2023-08-04) is inserted as row 387, pushing every existing row from
# 387 downward by one (old row 439 becomes row 440). This mirrors the
# commit "Fruta / hortaliza, semanal" which appends the latest week's
# reading to the top of this block of records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 387..439 down to 388..440, leaving a blank row 387 (keeps the
# existing date-format style on column D, same as Excel's native Insert).
$ws.Rows.Item(387).Insert()

# Populate the newly inserted row 387 with the new observation.
$ws.Cells.Item(387, 1).Value  = 4
$ws.Cells.Item(387, 2).Value  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(387, 3).Value  = 'Los Lagos'
$ws.Cells.Item(387, 4).Value  = 45142
$ws.Cells.Item(387, 5).Value  = 10
$ws.Cells.Item(387, 6).Value  = 100112032
$ws.Cells.Item(387, 7).Value  = 'Zapallo italiano'
$ws.Cells.Item(387, 8).Value  = 'Sin especificar'
$ws.Cells.Item(387, 9).Value  = 'Primera'
$ws.Cells.Item(387, 10).Value = 240
$ws.Cells.Item(387, 11).Value = 19000
$ws.Cells.Item(387, 12).Value = 20000
$ws.Cells.Item(387, 13).Value = 19500
$ws.Cells.Item(387, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(387, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(387, 16).Value = 390
$ws.Cells.Item(387, 17).Value = 50
$ws.Cells.Item(387, 18).Value = 'Hortaliza'
